$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so that numeric-looking
# strings (e.g. "536.60", "6.73") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.713.61"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.647.14"
$ws.Range("E3").Value = "  +2.06%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.16%  "

# Row 5 - BNB
$ws.Range("D5").Value = "536.60"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6 - Solana
$ws.Range("D6").Value = "145.28"
$ws.Range("E6").Value = "  +3.37%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.78%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.663.15"
$ws.Range("E9").Value = "  +2.14%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  +4.34%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.05%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.19%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.07%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.123.32"
$ws.Range("E14").Value = "  +2.34%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "59.623.98"
$ws.Range("E15").Value = "  +0.85%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "21.25"
$ws.Range("E16").Value = "  +3.91%  "

# Row 17 - was WrappedEther, now ShibaInu
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  +1.12%  "

# Row 18 - was ShibaInu, now WrappedEther
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.616.39"
$ws.Range("E18").Value = "  +2.68%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "344.92"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "4.41"
$ws.Range("E20").Value = "  +1.88%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  +1.21%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.35"
$ws.Range("E22").Value = "  -0.66%  "

# Row 23 - Dai
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.08%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "66.57"
$ws.Range("E24").Value = "  -1.38%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +2.24%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.90%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.43%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +1.98%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0749"
$ws.Range("E29").Value = "  +2.28%  "

# Row 30 - USDe
$ws.Range("E30").Value = "  -0.01%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.39%  "

# Row 32 - Aptos
$ws.Range("E32").Value = "  +1.03%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "19.03"
$ws.Range("E33").Value = "  +1.36%  "

# Row 34 - Monero
$ws.Range("D34").Value = "150.15"
$ws.Range("E34").Value = "  +0.25%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "4.03"
$ws.Range("E35").Value = "  +1.59%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +2.77%  "

# Row 37 - Stacks
$ws.Range("D37").Value = "1.46"
$ws.Range("E37").Value = "  -0.68%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "0.839"
$ws.Range("E38").Value = "  +0.49%  "

# Row 39 - was Bittensor, now SuiNetwork
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  -0.21%  "

# Row 40 - was SuiNetwork, now Bittensor
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "296.01"
$ws.Range("E40").Value = "  +8.84%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +1.84%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  -0.06%  "

# Row 43 - Mantle
$ws.Range("D43").Value = "0.604"
$ws.Range("E43").Value = "  +1.27%  "

# Row 44 - Hedera
$ws.Range("D44").Value = "0.0544"
$ws.Range("E44").Value = "  +4.91%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "19.40"
$ws.Range("E45").Value = "  +5.46%  "

# Row 46 - WhiteBITCoin
$ws.Range("E46").Value = "  -0.25%  "

# Row 47 - Stellar
$ws.Range("D47").Value = "0.0953"
$ws.Range("E47").Value = "  -0.70%  "

# Row 48 - VeChain
$ws.Range("D48").Value = "0.0226"
$ws.Range("E48").Value = "  +2.21%  "

# Row 49 - Maker
$ws.Range("D49").Value = "1.968.85"
$ws.Range("E49").Value = "  +1.12%  "

# Row 50 - RenderToken
$ws.Range("D50").Value = "4.55"
$ws.Range("E50").Value = "  +1.84%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "18.38"
$ws.Range("E51").Value = "  +0.69%  "

# Restore the original (default) style on the Price column now that all
# values have been written as text, so no visible formatting change remains.
$ws.Range("D2:D51").Style = "Normal"
